# Kenya_aggregated.xlsx: re-sort donor rows (A2:B52) into the new
# "Org type, flow date type" order. Each donor keeps its original
# amountUSD; only the row it lives on changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; Donor="United States of America, Government of"; Amount=187603597},
    @{Row=3; Donor="European Commission's Humanitarian Aid and Civil Protection Department"; Amount=39415742},
    @{Row=4; Donor="Qatar Charity"; Amount=533522},
    @{Row=5; Donor="Germany, Government of"; Amount=18190910},
    @{Row=6; Donor="New Zealand, Government of"; Amount=503597},
    @{Row=7; Donor="Canada, Government of"; Amount=10987333},
    @{Row=8; Donor="Denmark, Government of"; Amount=7758213},
    @{Row=9; Donor="Norway, Government of"; Amount=2833729},
    @{Row=10; Donor="Sweden, Government of"; Amount=9015988},
    @{Row=11; Donor="Ireland, Government of"; Amount=400082},
    @{Row=12; Donor="Private (individuals & organizations)"; Amount=1454955},
    @{Row=13; Donor="Central Emergency Response Fund"; Amount=16774818},
    @{Row=14; Donor="UNICEF National Committee/Sweden"; Amount=1238791},
    @{Row=15; Donor="UNICEF National Committee/Canada"; Amount=1621352},
    @{Row=16; Donor="UNICEF National Committee/Netherlands"; Amount=751367},
    @{Row=17; Donor="UNICEF National Committee/Australia"; Amount=70672},
    @{Row=18; Donor="UNICEF National Committee/France"; Amount=654268},
    @{Row=19; Donor="UNICEF National Committee/Germany"; Amount=4159623},
    @{Row=20; Donor="UNICEF National Committee/Luxembourg"; Amount=115898},
    @{Row=21; Donor="UNICEF National Committee/Norway"; Amount=66225},
    @{Row=22; Donor="UNICEF National Committee/Slovenia"; Amount=97102},
    @{Row=23; Donor="UNICEF National Committee/United Kingdom"; Amount=2187626},
    @{Row=24; Donor="US Fund for UNICEF"; Amount=429907},
    @{Row=25; Donor="Japan, Government of"; Amount=9516908},
    @{Row=26; Donor="United Kingdom, Government of"; Amount=18795515},
    @{Row=27; Donor="UNICEF National Committee/Denmark"; Amount=17379},
    @{Row=28; Donor="UNICEF National Committee/Spain"; Amount=52965},
    @{Row=29; Donor="UNICEF National Committee/Finland"; Amount=1492751},
    @{Row=30; Donor="Australia, Government of"; Amount=1849364},
    @{Row=31; Donor="UNICEF National Committee/Italy"; Amount=9800},
    @{Row=32; Donor="UNICEF National Committee/Portugal"; Amount=5717},
    @{Row=33; Donor="UNICEF National Committee/Belgium"; Amount=113766},
    @{Row=34; Donor="UNICEF National Committee/Switzerland"; Amount=9800},
    @{Row=35; Donor="UNICEF National Committee/Iceland"; Amount=58097},
    @{Row=36; Donor="UNICEF National Committee/Turkey"; Amount=69181},
    @{Row=37; Donor="China, Government of"; Amount=26770682},
    @{Row=38; Donor=""; Amount=2975673},
    @{Row=39; Donor="World Food Programme"; Amount=85358162},
    @{Row=40; Donor="Kenya, Government of"; Amount=6246},
    @{Row=41; Donor="Luxembourg, Government of"; Amount=559910},
    @{Row=42; Donor="Hungary, Government of"; Amount=42706},
    @{Row=43; Donor="France, Government of"; Amount=1390966},
    @{Row=44; Donor="Switzerland, Government of"; Amount=2743558},
    @{Row=45; Donor="Italy, Government of"; Amount=533618},
    @{Row=46; Donor="United Nations Population Fund"; Amount=46836},
    @{Row=47; Donor="Austria, Government of"; Amount=1119821},
    @{Row=48; Donor="Office for the Coordination of Humanitarian Affairs"; Amount=136274},
    @{Row=49; Donor="European Commission Directorate General for Development"; Amount=10449127},
    @{Row=50; Donor="European Commission"; Amount=2088106},
    @{Row=51; Donor="Belgium, Government of"; Amount=140499},
    @{Row=52; Donor="ACT Alliance / Church of Sweden"; Amount=974896}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Donor
    $ws.Cells.Item($item.Row, 2).Value = $item.Amount
}

Write-Output "Reordered $($data.Count) donor rows"
